# Apply the "Acme mAb" -> "COVIC" rebrand on the Dataset sheet, update the
# validation comments to match, and adjust the highlighted rows so that the
# example reflects a freshly-submitted (partially filled in) dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dataset")

# --- Rename antibodies from "Acme mAb N" to "COVIC N" ------------------
$ws.Range("A2").Value = "COVIC 1"
$ws.Range("A4").Value = "COVIC 1"
$ws.Range("A5").Value = "COVIC 4"
$ws.Range("A6").Value = "COVIC 5"
$ws.Range("A7").Value = "COVIC 6"
$ws.Range("A8").Value = "COVIC 7"
$ws.Range("A9").Value = "COVIC 8"
$ws.Range("A10").Value = "COVIC 9"
$ws.Range("A11").Value = "COVIC 10"

# --- Row 6: fill in the previously-missing measure with a typo value ---
# (keeps its existing highlighted style)
$ws.Range("B6").Value = "postive"

# --- Rows 8-11: clear out the qualitative measure values ----------------
# Row 8 also loses its highlighted fill (no longer flagged).
$ws.Range("B8").Clear()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()

# --- Update validation comments -----------------------------------------
[void]$ws.Range("A3").Comment.Text("Missing required value for 'Antibody name'")
[void]$ws.Range("A4").Comment.Text("Duplicate value 'COVIC 1' is not allowed for 'Antibody name'")
[void]$ws.Range("B6").Comment.Text("'postive' is not a recognized value for 'Qualitative measure'")
# B7's comment text is unchanged.
# B8's comment was already removed above by Range("B8").Clear().
